$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "Sumner Consolidation is in compliance with the inspection and collection requirement of  Paragraph 45"
#    -> "Sumner Consolidation appeared to be in compliance with the inspection and collection requirement of  Paragraph 45"
Replace-Text "Sumner Consolidation is in compliance with the inspection and collection requirement of  Paragraph 45" "Sumner Consolidation appeared to be in compliance with the inspection and collection requirement of  Paragraph 45"

# 2. "and remove trash from the buildings" -> "and remove waste from the buildings"
Replace-Text "and remove trash from the buildings" "and remove waste from the buildings"

# 3. "caretakers begin picking up trash each day" -> "caretakers begin picking up waste each day"
Replace-Text "caretakers begin picking up trash each day" "caretakers begin picking up waste each day"

# 4. "Sumner Consolidation is in compliance with the storage and removal" -> "Sumner Consolidation appeared to be in compliance with the storage and removal"
Replace-Text "Sumner Consolidation is in compliance with the storage and removal" "Sumner Consolidation appeared to be in compliance with the storage and removal"

# 5. "Bulk trash sits in a yard" -> "Bulk waste sits in a yard"
Replace-Text "Bulk trash sits in a yard" "Bulk waste sits in a yard"

# 6. "In terms of storage, in addition to disposing of" -> "In addition to disposing of"
Replace-Text "In terms of storage, in addition to disposing of" "In addition to disposing of"

# 7. Big paragraph rewrite: remove last sentence about pests/trash bins, change "leave their trash" -> "leave their waste",
#    and "drop-off sites" -> "drop sites" (both occurrences), within one contiguous block.
Replace-Text "Tenants are asked by management to leave their trash in the front of each building, either in trash cans or in exposed trash bags for pick up by caretakers if they choose not to use the chutes. Most tenants dispose of their trash using the drop-off sites. Waste is taken to one of four exterior compactors after being taken from the drop-off sites. All exterior compactors are in good shape and do not require maintenance at the time of reporting. When the trash is not removed from the premises, it is stored in a way that prevents pests (e.g., trash bins)." "Tenants are asked by management to leave their waste in the front of each building, either in trash cans or in exposed trash bags for pick up by caretakers if they choose not to use the chutes. Most tenants dispose of their trash using the drop sites. Waste is taken to one of four exterior compactors after being taken from the drop sites. All exterior compactors are in good shape and do not require maintenance at the time of reporting."

# 8. "Sumner has two bulk containers and 31 interior compactor rooms. Of the 31 interior compactor rooms, two were
#    inaccessible: 67 Marcus Garvey Boulevard due to pests and 987 Myrtle Avenue due to flooding. Further information
#    is needed to see what the current statues is of the interior compactors. Sumner disposes of approximately 100 – 200
#    compactor bags (40 lbs. Bags)."
#    -> "Sumner consolidation has two bulk containers. The consolidation disposes of approximately 100 – 200 compactor bags (40 lb. bags)."
Replace-Text "Sumner has two bulk containers and 31 interior compactor rooms. Of the 31 interior compactor rooms, two were inaccessible: 67 Marcus Garvey Boulevard due to pests and 987 Myrtle Avenue due to flooding. Further information is needed to see what the current statues is of the interior compactors. Sumner disposes of approximately 100 – 200 compactor bags (40 lbs. Bags)." "Sumner consolidation has two bulk containers. The consolidation disposes of approximately 100 – 200 compactor bags (40 lb. bags)."

# 9. "Sumner reports that, if necessary, they can take the trash from the developments" -> "Sumner consolidation reports that, if necessary, they can take the waste from the developments"
Replace-Text "Sumner reports that, if necessary, they can take the trash from the developments" "Sumner consolidation reports that, if necessary, they can take the waste from the developments"

# 10. "that resident outreach was the primary way to improve trash management. " -> "...improve waste management. "
Replace-Text "that resident outreach was the primary way to improve trash management. " "that resident outreach was the primary way to improve waste management. "
